$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build full data array for rows 2-10, columns A:T
$data = @(
    @("ECs", "Spp1", "Cd44", "ECs", 2, 0.6666666666666666, 209.0063303333334, 627.018991, 0.6751081226665357, 0.6751081226665357, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 54630.11705081173, 491671.0534573055, 0.4704610570330462, 0.4704610570330462),
    @("ECs", "Spp1", "Cd44", "FAPs", 2, 0.6666666666666666, 209.0063303333334, 627.018991, 0.6751081226665357, 0.6751081226665357, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 6688.083019045716, 60192.74717141144, 0.05759611687704159, 0.05759611687704159),
    @("ECs", "Spp1", "Cd44", "sCs", 2, 0.6666666666666666, 209.0063303333334, 627.018991, 0.6751081226665357, 0.6751081226665357, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 17075.61215302328, 153680.5093772095, 0.1470509487564478, 0.1470509487564478),
    @("FAPs", "Spp1", "Cd44", "ECs", 3, 1, 1.401741666666666, 4.205225, 0.004527744128790482, 0.004527744128790482, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 366.3875213868916, 3297.487692482024, 0.003155238719973302, 0.003155238719973302),
    @("FAPs", "Spp1", "Cd44", "FAPs", 3, 1, 1.401741666666666, 4.205225, 0.004527744128790482, 0.004527744128790482, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 44.85493153709999, 403.6943838339, 0.0003862795769678006, 0.0003862795769678006),
    @("FAPs", "Spp1", "Cd44", "sCs", 3, 1, 1.401741666666666, 4.205225, 0.004527744128790482, 0.004527744128790482, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 114.5209190580917, 1030.688271522825, 0.0009862258318493787, 0.000986225831849379),
    @("sCs", "Spp1", "Cd44", "ECs", 3, 1, 99.18134533333334, 297.544036, 0.3203641332046738, 0.3203641332046737, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 25924.04017703977, 233316.3615933579, 0.2232514225242004, 0.2232514225242003),
    @("sCs", "Spp1", "Cd44", "FAPs", 3, 1, 99.18134533333334, 297.544036, 0.3203641332046738, 0.3203641332046737, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 3173.746318937136, 28563.71687043423, 0.02733151837425395, 0.02733151837425395),
    @("sCs", "Spp1", "Cd44", "sCs", 3, 1, 99.18134533333334, 297.544036, 0.3203641332046738, 0.3203641332046737, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 8103.018616833562, 72927.16755150206, 0.06978119230621943, 0.06978119230621943)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $data[$i][$j]
    }
}
